# Adds a new "2022-Q3" quarter to the 酒鬼酒 (000799) holdings workbook:
#   1. Insert a "2022-Q3" summary row at the top of the "总计" sheet (pushing the
#      existing quarters down by one row).
#   2. Insert a new "2022-Q3" worksheet (fund-holder breakdown) right after "总计",
#      populated from a duplicate of the "2022-Q2" sheet (so it inherits the header
#      row, the 0-based row-index column, and all existing formatting).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift the 8 existing data rows down one row, then fill in
#    the new 2022-Q3 figures at row 2.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $next = $r + 1
    $summary.Range("B$next").Value = $summary.Range("B$r").Value()
    $summary.Range("C$next").Value = $summary.Range("C$r").Value()
    $summary.Range("D$next").Value = $summary.Range("D$r").Value()
}

# Row 9 needs the same index style as the rest of column A (A2:A8); A9 was
# previously unused so it has to be seeded from a neighbour before the value
# is overwritten.
$summary.Range("A8").Copy($summary.Range("A9"))
$summary.Range("A9").Value = 7

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 31
$summary.Range("D2").Value = 41.93

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate "2022-Q2" (After "总计") so the header
#    row / index column / styles come along for free, then overwrite the body
#    with the 2022-Q3 fund table and trim the sheet down to 32 rows.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$afterSheet = $wb.Worksheets.Item(1)
$template.Copy($null, $afterSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# row, col -> value ; [B..H] per data row (row 2 .. row 32). Column A already
# holds the correct 0-based index (0..30) copied from the "2022-Q2" template.
$dataRows = @(
    @("161725", "招商中证白酒指数A", "609.75", "94.82", "2.85", "17.3779", 8),
    @("180012", "银华富裕主题混合", "170.60", "94.25", "6.32", "10.7819", 8),
    @("512690", "鹏华中证酒ETF", "102.29", "97.72", "3.25", "3.3244", 9),
    @("012414", "招商中证白酒指数C", "92.42", "94.82", "2.85", "2.6340", 8),
    @("009542", "银华富利精选混合A", "28.27", "92.86", "6.84", "1.9337", 6),
    @("160632", "鹏华中证酒指数（LOF）A", "37.38", "94.22", "3.24", "1.2111", 8),
    @("010846", "南方卓越优选3个月持有期混合A", "20.44", "80.88", "5.86", "1.1978", 4),
    @("160222", "国泰国证食品饮料行业指数（LOF）A", "60.76", "91.77", "1.55", "0.9418", 10),
    @("012043", "鹏华中证酒指数（LOF）C", "23.02", "94.22", "3.24", "0.7458", 8),
    @("002851", "南方品质优选灵活配置混合A", "11.81", "71.72", "5.05", "0.5964", 4),
    @("487021", "工银优质精选混合", "8.54", "76.59", "2.50", "0.2135", 9),
    @("001140", "工银总回报灵活配置混合A", "6.10", "79.57", "3.01", "0.1836", 10),
    @("010847", "南方卓越优选3个月持有期混合C", "3.06", "80.88", "5.86", "0.1793", 4),
    @("001496", "工银聚焦30股票", "3.56", "84.46", "3.31", "0.1178", 9),
    @("000763", "工银新财富灵活配置混合", "2.76", "92.61", "3.19", "0.0880", 9),
    @("011903", "南方领航优选混合A", "1.56", "82.01", "5.37", "0.0838", 5),
    @("159758", "华夏中证红利质量ETF", "1.69", "99.23", "3.03", "0.0512", 10),
    @("004703", "南方兴盛先锋灵活配置混合", "0.97", "82.88", "5.22", "0.0506", 6),
    @("233008", "大摩消费领航混合基金", "0.95", "79.72", "5.14", "0.0488", 8),
    @("762001", "国金国鑫灵活配置混合A", "0.90", "90.18", "4.42", "0.0398", 8),
    @("009762", "国金国鑫灵活配置混合C", "0.73", "90.18", "4.42", "0.0323", 8),
    @("011904", "南方领航优选混合C", "0.48", "82.01", "5.37", "0.0258", 5),
    @("015040", "国泰国证食品饮料行业指数（LOF）C", "1.25", "91.77", "1.55", "0.0194", 10),
    @("159789", "建信中证饮料主题ETF", "0.58", "97.53", "3.12", "0.0181", 9),
    @("002159", "东吴国企改革主题灵活配置混合A", "0.17", "92.20", "8.07", "0.0137", 8),
    @("012615", "东吴国企改革主题灵活配置混合C", "0.13", "92.20", "8.07", "0.0105", 8),
    @("012763", "华泰紫金中证细分食品饮料产业主题指数A", "0.16", "95.00", "1.71", "0.0027", 10),
    @("014044", "银华富利精选混合C", "0.03", "92.86", "6.84", "0.0021", 6),
    @("011477", "工银总回报灵活配置混合C", "0.05", "79.57", "3.01", "0.0015", 10),
    @("012764", "华泰紫金中证细分食品饮料产业主题指数C", "0.08", "95.00", "1.71", "0.0014", 10),
    @("013501", "南方品质优选灵活配置混合C", "0.00", "71.72", "5.05", 0, 4)
)

$row = 2
foreach ($item in $dataRows) {
    $q3.Range("B$row").NumberFormat = "@"
    $q3.Range("B$row").Value = "'" + $item[0]
    $q3.Range("C$row").NumberFormat = "@"
    $q3.Range("C$row").Value = "'" + $item[1]
    $q3.Range("D$row").NumberFormat = "@"
    $q3.Range("D$row").Value = "'" + $item[2]
    $q3.Range("E$row").NumberFormat = "@"
    $q3.Range("E$row").Value = "'" + $item[3]
    $q3.Range("F$row").NumberFormat = "@"
    $q3.Range("F$row").Value = "'" + $item[4]

    if ($item[5] -is [string]) {
        $q3.Range("G$row").NumberFormat = "@"
        $q3.Range("G$row").Value = "'" + $item[5]
    } else {
        $q3.Range("G$row").Value = $item[5]
    }

    $q3.Range("H$row").Value = $item[6]
    $row++
}

# Drop the unused tail of the copied "2022-Q2" sheet (rows 33..110) so the
# dimensions match the 32-row 2022-Q3 table.
$q3.Range("A33:H110").Clear()

Write-Host "2022-Q3 quarter added"
